# Update the "types" worksheet:
#  - Change the type of several rows from "critical" to "regular"
#  - Remove the duplicated "regular" block of rows that followed (rows 17-31)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose type needs to change from "critical" to "regular"
$ws.Range("A2").Value = "regular"
$ws.Range("A3").Value = "regular"
$ws.Range("A6").Value = "regular"
$ws.Range("A8").Value = "regular"
$ws.Range("A9").Value = "regular"
$ws.Range("A10").Value = "regular"
$ws.Range("A11").Value = "regular"

# Remove the now-redundant "regular" rows (previously rows 17-31)
$ws.Rows("17:31").Delete()

# Match the new active selection
$ws.Range("A10").Select()
